$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report is regenerated: the row describing file
# "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md" moved from "Handed back" to a
# fresh "Ready for handoff" state (new handoff timestamps / handoff file
# refs), while the row describing "c0635734-...md" is unchanged in content
# but now sorts first. Concretely each of the three worksheets has its two
# data rows swapped (by identity), and the "048510b4" row picks up new
# status / datetime / handoff-file / error-detail values.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Hyperlinks must be rebuilt because changing Range.Value alone does not
# refresh the stored hyperlink "display" text in this engine.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md", [Type]::Missing, [Type]::Missing, "e2e\c0635734-0671-4fae-97fd-ad4939ff1ade.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md", [Type]::Missing, [Type]::Missing, "e2e\048510b4-dbe2-4a4a-8231-6b20506ae6b2.md")

$wsOverview.Range("A2").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.md"
$wsOverview.Range("A3").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 08:51:53"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md", [Type]::Missing, [Type]::Missing, "c0635734-0671-4fae-97fd-ad4939ff1ade.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f03a12db13e59218cca12fe503655e661e070834/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md", [Type]::Missing, [Type]::Missing, "c0635734-0671-4fae-97fd-ad4939ff1ade.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md", [Type]::Missing, [Type]::Missing, "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f03a12db13e59218cca12fe503655e661e070834/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md", [Type]::Missing, [Type]::Missing, "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md")

$wsZhCn.Range("G2").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.zh-cn.xlf"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 08:51:47"
$wsZhCn.Range("J3").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.zh-cn.xlf"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3e9d480f71b169c1f50867c59bee7e6b665fab3/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md."

# Error Detail column needs to be wide enough for the new long message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1667

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md", [Type]::Missing, [Type]::Missing, "c0635734-0671-4fae-97fd-ad4939ff1ade.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0083840a9a19e93e78c105eab41a1f8da80e3f73/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md", [Type]::Missing, [Type]::Missing, "c0635734-0671-4fae-97fd-ad4939ff1ade.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md", [Type]::Missing, [Type]::Missing, "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0083840a9a19e93e78c105eab41a1f8da80e3f73/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md", [Type]::Missing, [Type]::Missing, "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md")

$wsDeDe.Range("G2").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.de-de.xlf"
$wsDeDe.Range("J2").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.de-de.xlf"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 08:51:53"
$wsDeDe.Range("J3").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.de-de.xlf"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3e9d480f71b169c1f50867c59bee7e6b665fab3/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md."

$wsDeDe.Columns.Item(16).ColumnWidth = 39.1667
